# Scheduled data refresh: update market-board derived numeric columns
# (currentAveragePrice / NQ / HQ / LevePrice / LeveProfit columns) across
# several rows on multiple sheets, per the latest FFXIV marketboard snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3688.861
$ws.Range("I76").Value = 3503.9583
$ws.Range("K76").Value = 3503.9583
$ws.Range("M76").Value = -3188.9583

$ws.Range("H79").Value = 3688.861
$ws.Range("I79").Value = 3503.9583
$ws.Range("K79").Value = 3503.9583
$ws.Range("M79").Value = -2411.9583

$ws.Range("H95").Value = 60000
$ws.Range("J95").Value = 60000
$ws.Range("L95").Value = 60000
$ws.Range("N95").Value = -65492

$ws.Range("H116").Value = 3000
$ws.Range("I116").Value = 3000
$ws.Range("K116").Value = 3000
$ws.Range("M116").Value = 442

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7912.8306
$ws.Range("I32").Value = 7166.6772
$ws.Range("K32").Value = 7166.6772
$ws.Range("M32").Value = -6879.6772

$ws.Range("H63").Value = 2878.5
$ws.Range("I63").Value = 1742.7142
$ws.Range("K63").Value = 1742.7142
$ws.Range("M63").Value = -1056.7142

$ws.Range("H66").Value = 2878.5
$ws.Range("I66").Value = 1742.7142
$ws.Range("K66").Value = 8713.571
$ws.Range("M66").Value = -5281.571

$ws.Range("H95").Value = 33071.668
$ws.Range("J95").Value = 33071.668
$ws.Range("L95").Value = 33071.668
$ws.Range("N95").Value = -38563.668

$ws.Range("H96").Value = 19000
$ws.Range("J96").Value = 19000
$ws.Range("L96").Value = 19000
$ws.Range("N96").Value = -24492

$ws.Range("H110").Value = 1647
$ws.Range("I110").Value = 1355.5454
$ws.Range("J110").Value = 3250
$ws.Range("K110").Value = 1355.5454
$ws.Range("L110").Value = 3250
$ws.Range("M110").Value = 689.4546
$ws.Range("N110").Value = -7340

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 1149.2858
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

$ws.Range("H54").Value = 19663.5
$ws.Range("J54").Value = 27078.857
$ws.Range("L54").Value = 27078.857
$ws.Range("N54").Value = -28046.857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3446.111
$ws.Range("I62").Value = 3476.875
$ws.Range("J62").Value = 3200
$ws.Range("K62").Value = 3476.875
$ws.Range("L62").Value = 3200
$ws.Range("M62").Value = -2852.875
$ws.Range("N62").Value = -4448

$ws.Range("H65").Value = 3446.111
$ws.Range("I65").Value = 3476.875
$ws.Range("J65").Value = 3200
$ws.Range("K65").Value = 17384.375
$ws.Range("L65").Value = 16000
$ws.Range("M65").Value = -14264.375
$ws.Range("N65").Value = -22240

$ws.Range("H99").Value = 2487.6365
$ws.Range("I99").Value = 2556.25
$ws.Range("J99").Value = 2304.6667
$ws.Range("K99").Value = 2556.25
$ws.Range("L99").Value = 2304.6667
$ws.Range("M99").Value = -1058.25
$ws.Range("N99").Value = -5300.6667

$ws.Range("H105").Value = 468.9091
$ws.Range("I105").Value = 468.9091
$ws.Range("K105").Value = 468.9091
$ws.Range("M105").Value = 1278.0909

$ws.Range("H126").Value = 2487.6365
$ws.Range("I126").Value = 2556.25
$ws.Range("J126").Value = 2304.6667
$ws.Range("K126").Value = 7668.75
$ws.Range("L126").Value = 6914.000100000001
$ws.Range("M126").Value = -5198.75
$ws.Range("N126").Value = -11854.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 99.117645
$ws.Range("I38").Value = 38
$ws.Range("K38").Value = 114
$ws.Range("M38").Value = 233

$ws.Range("H63").Value = 3647.5
$ws.Range("I63").Value = 3395
$ws.Range("J63").Value = 3900
$ws.Range("K63").Value = 10185
$ws.Range("L63").Value = 11700
$ws.Range("M63").Value = -9436
$ws.Range("N63").Value = -13198

$ws.Range("H66").Value = 3647.5
$ws.Range("I66").Value = 3395
$ws.Range("J66").Value = 3900
$ws.Range("K66").Value = 30555
$ws.Range("L66").Value = 35100
$ws.Range("M66").Value = -26811
$ws.Range("N66").Value = -42588

$ws.Range("H113").Value = 730.4194
$ws.Range("I113").Value = 745.42645
$ws.Range("J113").Value = 689.6
$ws.Range("K113").Value = 2236.27935
$ws.Range("L113").Value = 2068.8
$ws.Range("M113").Value = -66.27935000000025
$ws.Range("N113").Value = -6408.8

$ws.Range("H129").Value = 1840.6923
$ws.Range("I129").Value = 2616
$ws.Range("J129").Value = 1356.125
$ws.Range("K129").Value = 7848
$ws.Range("L129").Value = 4068.375
$ws.Range("M129").Value = -2848
$ws.Range("N129").Value = -14068.375

$ws.Range("H137").Value = 38466316
$ws.Range("I137").Value = 55558276
$ws.Range("J137").Value = 9397
$ws.Range("K137").Value = 166674828
$ws.Range("L137").Value = 28191
$ws.Range("M137").Value = -166669728
$ws.Range("N137").Value = -38391

$ws.Range("H138").Value = 5117.294
$ws.Range("I138").Value = 9410.833000000001
$ws.Range("J138").Value = 2775.3635
$ws.Range("K138").Value = 28232.499
$ws.Range("L138").Value = 8326.0905
$ws.Range("M138").Value = -23092.499
$ws.Range("N138").Value = -18606.0905

$ws.Range("H140").Value = 1821.625
$ws.Range("I140").Value = 1515.1904
$ws.Range("K140").Value = 4545.5712
$ws.Range("M140").Value = 634.4287999999997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5266.1494
$ws.Range("I70").Value = 4782.4473
$ws.Range("K70").Value = 4782.4473
$ws.Range("M70").Value = -4512.4473

$ws.Range("H73").Value = 5266.1494
$ws.Range("I73").Value = 4782.4473
$ws.Range("K73").Value = 4782.4473
$ws.Range("M73").Value = -3846.4473

$ws.Range("H102").Value = 3729.5366
$ws.Range("I102").Value = 3411.5518
$ws.Range("J102").Value = 4498
$ws.Range("K102").Value = 3411.5518
$ws.Range("L102").Value = 4498
$ws.Range("M102").Value = -1789.5518
$ws.Range("N102").Value = -7742

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4526.2607
$ws.Range("I40").Value = 4228
$ws.Range("J40").Value = 5600
$ws.Range("K40").Value = 4228
$ws.Range("L40").Value = 5600
$ws.Range("M40").Value = -4092
$ws.Range("N40").Value = -5872

$ws.Range("H70").Value = 30163
$ws.Range("J70").Value = 30163
$ws.Range("L70").Value = 30163
$ws.Range("N70").Value = -30703

$ws.Range("H73").Value = 30163
$ws.Range("J73").Value = 30163
$ws.Range("L73").Value = 30163
$ws.Range("N73").Value = -32035

$ws.Range("H122").Value = 7181.4814
$ws.Range("I122").Value = 5616.6665
$ws.Range("K122").Value = 16849.9995
$ws.Range("M122").Value = -14399.9995

$ws.Range("H132").Value = 3004.0322
$ws.Range("I132").Value = 2632.5557
$ws.Range("J132").Value = 3518.3845
$ws.Range("K132").Value = 7897.6671
$ws.Range("L132").Value = 10555.1535
$ws.Range("M132").Value = -5367.6671
$ws.Range("N132").Value = -15615.1535

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 100049
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()

$ws.Range("H132").Value = 2608.239
$ws.Range("I132").Value = 2268.5938
$ws.Range("J132").Value = 3384.5715
$ws.Range("K132").Value = 6805.7814
$ws.Range("L132").Value = 10153.7145
$ws.Range("M132").Value = -4275.7814
$ws.Range("N132").Value = -15213.7145
